$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H) - copy formatting from the adjacent "sum" header (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the Save values for each data row
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
